$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 307.44446
$ws.Range("I2").Value = 309.7143
$ws.Range("J2").Value = 299.5
$ws.Range("K2").Value = 309.7143
$ws.Range("L2").Value = 299.5
$ws.Range("M2").Value = -196.7143
$ws.Range("N2").Value = -525.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 74
$ws.Range("I5").Value = 76.416664
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 76.416664
$ws.Range("L5").Value = 45
$ws.Range("M5").Value = 38.583336
$ws.Range("N5").Value = -275

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 210.5
$ws.Range("I12").Value = 210.5
$ws.Range("K12").Value = 210.5
$ws.Range("M12").Value = -40.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 177.85715
$ws.Range("I55").Value = 189.2
$ws.Range("K55").Value = 189.2
$ws.Range("M55").Value = 24.80000000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1166.6666
$ws.Range("J125").Value = 1200
$ws.Range("L125").Value = 10800
$ws.Range("N125").Value = -15720

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2881.3333
$ws.Range("I127").Value = 1857.6
$ws.Range("K127").Value = 5572.799999999999
$ws.Range("M127").Value = -612.7999999999993

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 6334.125
$ws.Range("I131").Value = 2670
$ws.Range("K131").Value = 8010
$ws.Range("M131").Value = -2970

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 100000
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 100000
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 522.5
$ws.Range("I12").Value = 105
$ws.Range("K12").Value = 105
$ws.Range("M12").Value = 63

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2859.2
$ws.Range("I99").Value = 1498.6666
$ws.Range("J99").Value = 4900
$ws.Range("K99").Value = 1498.6666
$ws.Range("L99").Value = 4900
$ws.Range("M99").Value = -0.6666000000000167
$ws.Range("N99").Value = -7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 359.875
$ws.Range("I5").Value = 359.875
$ws.Range("K5").Value = 359.875
$ws.Range("M5").Value = -247.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 11500.5
$ws.Range("I8").Value = 10001
$ws.Range("K8").Value = 10001
$ws.Range("M8").Value = -9861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2063
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 270.8
$ws.Range("I12").Value = 336
$ws.Range("K12").Value = 336
$ws.Range("M12").Value = -166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 399.25
$ws.Range("J15").Value = 432.33334
$ws.Range("L15").Value = 432.33334
$ws.Range("N15").Value = -772.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1200.6
$ws.Range("I107").Value = 1158.5264
$ws.Range("K107").Value = 1158.5264
$ws.Range("M107").Value = 761.4736

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 303425.78
$ws.Range("J141").Value = 303425.78
$ws.Range("L141").Value = 303425.78
$ws.Range("N141").Value = -313785.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2525
$ws.Range("I2").Value = 50
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7537683
$ws.Range("I4").Value = 4328780.5
$ws.Range("K4").Value = 12986341.5
$ws.Range("M4").Value = -12986229.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78.71429000000001
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 840.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 400.42856
$ws.Range("I50").Value = 375.5
$ws.Range("K50").Value = 1126.5
$ws.Range("M50").Value = -645.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 400.42856
$ws.Range("I53").Value = 375.5
$ws.Range("K53").Value = 1126.5
$ws.Range("M53").Value = -645.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 979.125
$ws.Range("I108").Value = 1258.8334
$ws.Range("J108").Value = 140
$ws.Range("K108").Value = 3776.5002
$ws.Range("L108").Value = 420
$ws.Range("M108").Value = -896.5001999999999
$ws.Range("N108").Value = -6180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1475.5834
$ws.Range("I113").Value = 1787
$ws.Range("J113").Value = 1039.6
$ws.Range("K113").Value = 5361
$ws.Range("L113").Value = 3118.8
$ws.Range("M113").Value = -3191
$ws.Range("N113").Value = -7458.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3877.7778
$ws.Range("I121").Value = 633.3333
$ws.Range("K121").Value = 1899.9999
$ws.Range("M121").Value = -589.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 9847.5
$ws.Range("J6").Value = 9847.5
$ws.Range("L6").Value = 9847.5
$ws.Range("N6").Value = -10073.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H16").Value = 9847.5
$ws.Range("J16").Value = 9847.5
$ws.Range("L16").Value = 9847.5
$ws.Range("N16").Value = -10347.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 6998
$ws.Range("I22").Value = 6999
$ws.Range("J22").Value = 6997.5
$ws.Range("K22").Value = 6999
$ws.Range("L22").Value = 6997.5
$ws.Range("M22").Value = -6470
$ws.Range("N22").Value = -8055.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 5299.5713
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 23996
$ws.Range("J95").Value = 23996
$ws.Range("L95").Value = 23996
$ws.Range("N95").Value = -29488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2514.16
$ws.Range("I102").Value = 1811.591
$ws.Range("K102").Value = 1811.591
$ws.Range("M102").Value = -189.5909999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2177.375
$ws.Range("I9").Value = 488.42856
$ws.Range("K9").Value = 488.42856
$ws.Range("M9").Value = -264.42856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 15000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -15590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1508
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1252

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 9250
$ws.Range("J39").Value = 11666.667
$ws.Range("L39").Value = 11666.667
$ws.Range("N39").Value = -12586.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 601.6667
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 5500
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
